$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric,
# so Excel keeps them as literal text (preserving trailing zeros / exact digits)
# instead of auto-converting them to a number.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D18", "D19", "D21", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range('D2').Value = '26.610.12'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '1.827.21'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '308.80'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').Value = '0.4689'
$ws.Range('E7').Value = '  +3.39%  '
$ws.Range('D8').Value = '0.3596'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').Value = '0.07139'
$ws.Range('E9').Value = '  +0.50%  '
$ws.Range('D10').Value = '0.9286'
$ws.Range('E10').Value = '  +4.28%  '
$ws.Range('D11').Value = '0.07660'
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').Value = '19.43'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('D13').Value = '1.834.60'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').Value = '5.255'
$ws.Range('E14').Value = '  -0.86%  '
$ws.Range('D15').Value = '6.333'
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('D16').Value = '87.41'
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = '0.000008537'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').Value = '1.007'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').Value = '26.630.84'
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('D21').Value = '14.24'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('E22').Value = '  +0.65%  '
$ws.Range('D23').Value = '2.076.90'
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('D24').Value = '10.56'
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = '1.912'
$ws.Range('E25').Value = '  -2.97%  '
$ws.Range('D26').Value = '151.84'
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('D28').Value = '1.991'
$ws.Range('E28').Value = '  -3.42%  '
$ws.Range('D29').Value = '113.44'
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('D30').Value = '4.856'
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('D31').Value = '0.08819'
$ws.Range('E31').Value = '  +1.19%  '
$ws.Range('D32').Value = '3.158'
$ws.Range('E32').Value = '  +0.93%  '
$ws.Range('D33').Value = '2.826'
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('D34').Value = '1.159'
$ws.Range('E34').Value = '  +4.28%  '
$ws.Range('D35').Value = '0.7371'
$ws.Range('E35').Value = '  +1.76%  '
$ws.Range('D36').Value = '4.434'
$ws.Range('E36').Value = '  -0.31%  '
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').Value = '0.01921'
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('D39').Value = '2.939'
$ws.Range('E39').Value = '  +1.60%  '
$ws.Range('D40').Value = '0.05146'
$ws.Range('E40').Value = '  +0.79%  '
$ws.Range('D41').Value = '6.899'
$ws.Range('E41').Value = '  +1.35%  '
$ws.Range('D42').Value = '0.5052'
$ws.Range('E42').Value = '  -2.37%  '
$ws.Range('D43').Value = '0.1495'
$ws.Range('E43').Value = '  -1.35%  '
$ws.Range('D44').Value = '8.075'
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('D46').Value = '0.4630'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').Value = '10.09'
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('D48').Value = '98.38'
$ws.Range('E48').Value = '  -2.94%  '
$ws.Range('D49').Value = '1.572'
$ws.Range('E49').Value = '  -0.31%  '
$ws.Range('D50').Value = '0.06024'
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('D51').Value = '63.78'
$ws.Range('E51').Value = '  -0.51%  '
